# Auto-generated edit script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '47.226.78'
$ws.Range("E2").Value = '  -0.17%  '

# Row 3
$ws.Range("D3").Value = '2.487.05'
$ws.Range("E3").Value = '  -0.86%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = "'321.00"
$ws.Range("E5").Value = '  -0.98%  '

# Row 6
$ws.Range("D6").Value = "'107.71"
$ws.Range("E6").Value = '  +1.61%  '

# Row 7
$ws.Range("D7").Value = "'0.521"
$ws.Range("E7").Value = '  -0.60%  '

# Row 8
$ws.Range("D8").Value = "'0.999"

# Row 9
$ws.Range("E9").Value = '  -1.55%  '

# Row 10
$ws.Range("D10").Value = "'38.47"
$ws.Range("E10").Value = '  +4.76%  '

# Row 11
$ws.Range("D11").Value = "'0.0808"
$ws.Range("E11").Value = '  -1.23%  '

# Row 12
$ws.Range("E12").Value = '  +0.08%  '

# Row 13
$ws.Range("E13").Value = '  -0.91%  '

# Row 14
$ws.Range("D14").Value = "'7.10"
$ws.Range("E14").Value = '  -1.06%  '

# Row 15
$ws.Range("D15").Value = '2.873.65'
$ws.Range("E15").Value = '  -0.96%  '

# Row 16
$ws.Range("D16").Value = '2.489.11'
$ws.Range("E16").Value = '  -0.45%  '

# Row 17
$ws.Range("D17").Value = "'0.845"
$ws.Range("E17").Value = '  -0.19%  '

# Row 18
$ws.Range("D18").Value = '47.141.64'
$ws.Range("E18").Value = '  -0.18%  '

# Row 19
$ws.Range("D19").Value = "'12.70"
$ws.Range("E19").Value = '  -0.92%  '

# Row 20
$ws.Range("E20").Value = '  +1.31%  '

# Row 21
$ws.Range("E21").Value = '  -1.29%  '

# Row 22
$ws.Range("E22").Value = '  +13.29%  '

# Row 23
$ws.Range("D23").Value = "'70.22"
$ws.Range("E23").Value = '  -1.02%  '

# Row 24
$ws.Range("D24").Value = "'245.03"
$ws.Range("E24").Value = '  -2.95%  '

# Row 25
$ws.Range("E25").Value = '  -0.08%  '

# Row 26
$ws.Range("E26").Value = '  -0.01%  '

# Row 27
$ws.Range("D27").Value = "'25.65"
$ws.Range("E27").Value = '  -2.89%  '

# Row 28
$ws.Range("E28").Value = '  -0.89%  '

# Row 29
$ws.Range("D29").Value = "'9.98"
$ws.Range("E29").Value = '  -0.05%  '

# Row 30
$ws.Range("E30").Value = '  -2.53%  '

# Row 31
$ws.Range("E31").Value = '  -0.49%  '

# Row 32
$ws.Range("D32").Value = "'49.45"
$ws.Range("E32").Value = '  -0.72%  '

# Row 33
$ws.Range("D33").Value = "'20.28"
$ws.Range("E33").Value = '  +2.01%  '

# Row 34
$ws.Range("E34").Value = '  -0.09%  '

# Row 35
$ws.Range("D35").Value = "'0.0777"
$ws.Range("E35").Value = '  -0.18%  '

# Row 36
$ws.Range("E36").Value = '  +0.07%  '

# Row 37
$ws.Range("E37").Value = '  +0.27%  '

# Row 38
$ws.Range("D38").Value = "'4.61"
$ws.Range("E38").Value = '  -0.69%  '

# Row 39
$ws.Range("D39").Value = "'2.92"
$ws.Range("E39").Value = '  -1.83%  '

# Row 40
$ws.Range("B40").Value = 'EnergySwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D40").Value = "'22.65"
$ws.Range("E40").Value = '  +3.88%  '

# Row 41
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").Value = "'0.111"
$ws.Range("E41").Value = '  -0.75%  '

# Row 42
$ws.Range("E42").Value = '  -0.66%  '

# Row 43
$ws.Range("D43").Value = "'118.87"
$ws.Range("E43").Value = '  -3.86%  '

# Row 44
$ws.Range("E44").Value = '  -1.14%  '

# Row 45
$ws.Range("D45").Value = '1.982.61'
$ws.Range("E45").Value = '  +0.17%  '

# Row 46
$ws.Range("D46").Value = "'3.00"
$ws.Range("E46").Value = '  -0.74%  '

# Row 47
$ws.Range("D47").Value = "'1.98"
$ws.Range("E47").Value = '  -6.88%  '

# Row 48
$ws.Range("E48").Value = '  -0.82%  '

# Row 49
$ws.Range("E49").Value = '  -2.73%  '

# Row 50
$ws.Range("D50").Value = "'5.12"
$ws.Range("E50").Value = '  -5.44%  '

# Row 51
$ws.Range("E51").Value = '  +3.12%  '
